$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.017.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.056.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.84"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.111"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.891"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.355.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.043.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.990.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.83"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.14"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +12.61%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0615"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0840"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.54%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.27"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0958"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -11.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.300.67"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.78"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.241.80"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.96%  "
